# Weekly update: insert a new price record as row 33 (pushing the
# existing rows 33-43 down to 34-44), matching the new week's data
# (Feria Lagunitas de Puerto Montt - Espinaca, fecha 2022-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 33:43 down to 34:44, duplicating row 33's formatting
# (including the date style on column D) into the newly freed row 33.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with this week's record.
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44813
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = 100112012
$ws.Cells.Item(33, 7).Value = "Espinaca"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 35
$ws.Cells.Item(33, 11).Value = 12000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 12000
$ws.Cells.Item(33, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(33, 15).Value = "Región Metropolitana"
$ws.Cells.Item(33, 16).Value = 1200
$ws.Cells.Item(33, 17).Value = 10
$ws.Cells.Item(33, 18).Value = "Hortaliza"
